$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("REGISTRO")
$ws2 = $wb.Worksheets.Item("PARAMETROS")

# --- New rows of activity log data in REGISTRO ---

# Row 5 - SE SIGUE EN CONFECCION DE ANTEPROYECTO (no time range given)
$ws1.Range("A5").Value = 44301
$ws1.Range("B5").Value = "DOCUMENTACION"
$ws1.Range("C5").Value = "SE SIGUE EN CONFECCION DE ANTEPROYECTO"

# Row 6 - ENTREVISTA CON SUBJEFE OPERATIVO SILVIA DELGADO
$ws1.Range("A6").Value = 44302
$ws1.Range("B6").Value = "REUNION"
$ws1.Range("C6").Value = "ENTREVISTA CON SUBJEFE OPERATIVO SILVIA DELGADO"
$ws1.Range("E6").Value = 0.54166666666666663
$ws1.Range("F6").Value = 0.625

# Row 7 - ENTREVISTA CON SUBJEFE OPERATIVO MARIEL FERREIRA
$ws1.Range("A7").Value = 44306
$ws1.Range("B7").Value = "REUNION"
$ws1.Range("C7").Value = "ENTREVISTA CON SUBJEFE OPERATIVO MARIEL FERREIRA"
$ws1.Range("E7").Value = 0.5625
$ws1.Range("F7").Value = 0.60416666666666663

# --- Formats ---

# Column A (dates) now uses a dd/mm/yyyy custom format
$ws1.Columns.Item(1).ColumnWidth = 10.6
$ws1.Range("A1:A7").NumberFormat = "dd/mm/yyyy;@"

# Time columns keep their h:mm format (re-applied to the new rows too)
$ws1.Range("E1:F4").NumberFormat = "h:mm"
$ws1.Range("E6:F7").NumberFormat = "h:mm"

# Page setup (A4/Letter-ish printer paper, portrait)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Sheet selection / activation ---
# PARAMETROS no longer the active tab, but its own selection moves to A5
$ws2.Range("A5").Select()
# REGISTRO becomes the active tab with its selection at E8
$ws1.Range("E8").Select()

Write-Output "done"
